$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up the "Ford" column values (drop the floating-point noise) ---
$ws.Range("E4").Value = 303333
$ws.Range("E5").Value = 453866
$ws.Range("E6").Value = 373206

# --- 2. Remove the calculated "Total" column from the table entirely ---
$tbl = $ws.ListObjects.Item("march")
$totalCol = $tbl.ListColumns.Item("Total")
$totalRange = $totalCol.Range
$totalCol.Delete()
$totalRange.EntireColumn.Delete()

# --- 3. Give the remaining numeric columns a 2-decimal number format ---
$dataRange = $ws.Range("B4:E7")
$dataRange.NumberFormat = "#,##0.00_ ;-#,##0.00 "

$ws.Range("E12").Select() | Out-Null
